# Commit: "Added prueba 4/8 testMansion"
# Mark the implemented tests (column D, "Test implementado") with an "x"
# for the "testMansion" suite rows (Anadir habitacion, Eliminar habitacion,
# Anadir pasillo, Eliminar Pasillo) on the "Conexion a interfaz" sheet, and
# make that sheet the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conexion a interfaz")

# Switch focus to this sheet (becomes the workbook's active tab / tabSelected).
$ws.Activate()

# Mark the 4/8 implemented tests.
$ws.Range("D5").Value = "x"
$ws.Range("D6").Value = "x"
$ws.Range("D9").Value = "x"
$ws.Range("D10").Value = "x"

# Leave the selection where the author left it.
$ws.Range("E5").Select()
